$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to match statsmodel's Tukey test output naming:
#   A    -> group1
#   B    -> group2
#   pval -> p-adj
$ws.Range("A1").Value = "group1"
$ws.Range("B1").Value = "group2"
$ws.Range("D1").Value = "p-adj"

# Update the active selection to D2 (matches the post-edit state)
$ws.Range("D2").Select()
